$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-36
# from serial date 45685 (2025-01-28) to 45686 (2025-01-29)
$ws.Range("C2:C36").Value = 45686
